$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 4016.5386
$ws.Cells.Item(76, 9).Value = 3049.8333
$ws.Cells.Item(76, 11).Value = 3049.8333
$ws.Cells.Item(76, 13).Value = -2734.8333
$ws.Cells.Item(79, 8).Value = 4016.5386
$ws.Cells.Item(79, 9).Value = 3049.8333
$ws.Cells.Item(79, 11).Value = 3049.8333
$ws.Cells.Item(79, 13).Value = -1957.8333
$ws.Cells.Item(86, 8).Value = 16999.8
$ws.Cells.Item(86, 9).Value = 16249.75
$ws.Cells.Item(86, 11).Value = 16249.75
$ws.Cells.Item(86, 13).Value = -15126.75
$ws.Cells.Item(89, 8).Value = 16999.8
$ws.Cells.Item(89, 9).Value = 16249.75
$ws.Cells.Item(89, 11).Value = 81248.75
$ws.Cells.Item(89, 13).Value = -75632.75
$ws.Cells.Item(99, 8).Value = 155.33333
$ws.Cells.Item(99, 9).Value = 155.33333
$ws.Cells.Item(99, 11).Value = 465.99999
$ws.Cells.Item(99, 13).Value = 1032.00001
$ws.Cells.Item(100, 8).Value = 2191.8572
$ws.Cells.Item(100, 9).Value = 1903.2727
$ws.Cells.Item(100, 10).Value = 3250
$ws.Cells.Item(100, 11).Value = 1903.2727
$ws.Cells.Item(100, 12).Value = 3250
$ws.Cells.Item(100, 13).Value = -1362.2727
$ws.Cells.Item(100, 14).Value = -4332
$ws.Cells.Item(101, 8).Value = 242
$ws.Cells.Item(101, 9).Value = 222.66667
$ws.Cells.Item(101, 11).Value = 668.00001
$ws.Cells.Item(101, 13).Value = 953.99999
$ws.Cells.Item(106, 8).Value = 5373.6665
$ws.Cells.Item(106, 9).Value = 5373.6665
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 5373.6665
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = -4742.6665
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 455.7
$ws.Cells.Item(115, 9).Value = 483.66666
$ws.Cells.Item(115, 11).Value = 1450.99998
$ws.Cells.Item(115, 13).Value = 116.0000199999999
$ws.Cells.Item(132, 8).Value = 11506.2
$ws.Cells.Item(132, 9).Value = 13632.6875
$ws.Cells.Item(132, 11).Value = 40898.0625
$ws.Cells.Item(132, 13).Value = -38368.0625
$ws.Cells.Item(138, 8).Value = 4565.3228
$ws.Cells.Item(138, 10).Value = 4855.8125
$ws.Cells.Item(138, 12).Value = 14567.4375
$ws.Cells.Item(138, 14).Value = -24847.4375
$ws.Cells.Item(141, 8).Value = 7156.5884
$ws.Cells.Item(141, 9).Value = 3026.2
$ws.Cells.Item(141, 10).Value = 13057.143
$ws.Cells.Item(141, 11).Value = 9078.599999999999
$ws.Cells.Item(141, 12).Value = 39171.429
$ws.Cells.Item(141, 13).Value = -3898.599999999999
$ws.Cells.Item(141, 14).Value = -49531.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 8344624.5
$ws.Cells.Item(61, 10).Value = 20003100
$ws.Cells.Item(61, 12).Value = 20003100
$ws.Cells.Item(61, 14).Value = -20003524
$ws.Cells.Item(97, 8).Value = 55557020
$ws.Cells.Item(97, 9).Value = 1431.9286
$ws.Cells.Item(97, 10).Value = 250001570
$ws.Cells.Item(97, 11).Value = 1431.9286
$ws.Cells.Item(97, 12).Value = 250001570
$ws.Cells.Item(97, 13).Value = -935.9286
$ws.Cells.Item(97, 14).Value = -250002562
$ws.Cells.Item(122, 8).Value = 2199.7917
$ws.Cells.Item(122, 9).Value = 1852.2
$ws.Cells.Item(122, 11).Value = 5556.6
$ws.Cells.Item(122, 13).Value = -3106.6
$ws.Cells.Item(136, 8).Value = 8344624.5
$ws.Cells.Item(136, 10).Value = 20003100
$ws.Cells.Item(136, 12).Value = 60009300
$ws.Cells.Item(136, 14).Value = -60014400

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 343.41177
$ws.Cells.Item(80, 9).Value = 451
$ws.Cells.Item(80, 11).Value = 451
$ws.Cells.Item(80, 13).Value = 547
$ws.Cells.Item(83, 8).Value = 343.41177
$ws.Cells.Item(83, 9).Value = 451
$ws.Cells.Item(83, 11).Value = 2255
$ws.Cells.Item(83, 13).Value = 2737
$ws.Cells.Item(105, 8).Value = 14816.0625
$ws.Cells.Item(105, 9).Value = 12223.143
$ws.Cells.Item(105, 10).Value = 16832.777
$ws.Cells.Item(105, 11).Value = 12223.143
$ws.Cells.Item(105, 12).Value = 16832.777
$ws.Cells.Item(105, 13).Value = -10476.143
$ws.Cells.Item(105, 14).Value = -20326.777
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 8256.357
$ws.Cells.Item(131, 9).Value = 2998
$ws.Cells.Item(131, 11).Value = 8994
$ws.Cells.Item(131, 13).Value = -3954

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(63, 8).Value = 50001
$ws.Cells.Item(63, 10).Value = 50001
$ws.Cells.Item(63, 12).Value = 50001
$ws.Cells.Item(63, 14).Value = -51373
$ws.Cells.Item(64, 8).Value = 35000
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 35000
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 35000
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -35496
$ws.Cells.Item(66, 8).Value = 50001
$ws.Cells.Item(66, 10).Value = 50001
$ws.Cells.Item(66, 12).Value = 150003
$ws.Cells.Item(66, 14).Value = -156867
$ws.Cells.Item(67, 8).Value = 35000
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 35000
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 35000
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -36716
$ws.Cells.Item(74, 8).Value = 49000
$ws.Cells.Item(74, 10).Value = 49000
$ws.Cells.Item(74, 12).Value = 49000
$ws.Cells.Item(74, 14).Value = -50872
$ws.Cells.Item(75, 8).Value = 50001
$ws.Cells.Item(75, 10).Value = 50001
$ws.Cells.Item(75, 12).Value = 50001
$ws.Cells.Item(75, 14).Value = -51749
$ws.Cells.Item(77, 8).Value = 49000
$ws.Cells.Item(77, 10).Value = 49000
$ws.Cells.Item(77, 12).Value = 147000
$ws.Cells.Item(77, 14).Value = -156360
$ws.Cells.Item(78, 8).Value = 50001
$ws.Cells.Item(78, 10).Value = 50001
$ws.Cells.Item(78, 12).Value = 150003
$ws.Cells.Item(78, 14).Value = -158739
$ws.Cells.Item(80, 8).Value = 1824.6666
$ws.Cells.Item(80, 9).Value = 1138.25
$ws.Cells.Item(80, 10).Value = 3197.5
$ws.Cells.Item(80, 11).Value = 1138.25
$ws.Cells.Item(80, 12).Value = 3197.5
$ws.Cells.Item(80, 13).Value = -140.25
$ws.Cells.Item(80, 14).Value = -5193.5
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 1824.6666
$ws.Cells.Item(83, 9).Value = 1138.25
$ws.Cells.Item(83, 10).Value = 3197.5
$ws.Cells.Item(83, 11).Value = 5691.25
$ws.Cells.Item(83, 12).Value = 15987.5
$ws.Cells.Item(83, 13).Value = -699.25
$ws.Cells.Item(83, 14).Value = -25971.5
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 2649.5
$ws.Cells.Item(32, 9).Value = 2649.5
$ws.Cells.Item(32, 11).Value = 2649.5
$ws.Cells.Item(32, 13).Value = -2332.5
$ws.Cells.Item(46, 8).Value = 2770.3333
$ws.Cells.Item(46, 9).Value = 1192.3846
$ws.Cells.Item(46, 10).Value = 4235.5713
$ws.Cells.Item(46, 11).Value = 1192.3846
$ws.Cells.Item(46, 12).Value = 4235.5713
$ws.Cells.Item(46, 13).Value = -1004.3846
$ws.Cells.Item(46, 14).Value = -4611.5713
$ws.Cells.Item(122, 8).Value = 5470.1665
$ws.Cells.Item(122, 9).Value = 3663.8572
$ws.Cells.Item(122, 11).Value = 10991.5716
$ws.Cells.Item(122, 13).Value = -8541.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 29959.5
$ws.Cells.Item(124, 10).Value = 29959.5
$ws.Cells.Item(124, 12).Value = 29959.5
$ws.Cells.Item(124, 14).Value = -39779.5
